# Update column F (dSF) values on the active sheet to reflect the
# repulled / recalculated data from the commit "repull data, push all
# data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    4  = 2
    6  = 1
    7  = 0
    9  = -1
    11 = 2
    12 = 1
    13 = 1
    14 = 2
    15 = -2
    16 = -1
    17 = -4
    19 = -2
    20 = -3
    21 = 2
    22 = -2
    23 = -1
    25 = -4
    26 = -5
    27 = -3
    28 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
